$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Default workbook font Calibri -> Tahoma (applies to the "Normal" style,
#    which backs every unstyled cell in the workbook).
# ---------------------------------------------------------------------------
$wb.Styles.Item("Normal").Font.Name = "Tahoma"

# ---------------------------------------------------------------------------
# 2) dcb2.0 sheet: collapse the four "DCB_xx to MPx" rows into a single
#    "DCB2.0" label row (matches the other single-label sheets).
#    Do this before touching the 5dh/nissan/custom sheets so the new
#    "DCB2.0" shared string is appended right after "DCB1.2H".
# ---------------------------------------------------------------------------
$wsDcb20 = $wb.Worksheets.Item("dcb2.0")
$wsDcb20.Rows.Item(3).Delete()
$wsDcb20.Rows.Item(3).Delete()
$wsDcb20.Rows.Item(3).Delete()
$wsDcb20.Range("A2").Value = "DCB2.0"
$wsDcb20.Range("B2").ClearContents()

# ---------------------------------------------------------------------------
# 3) bmw / dcb1.2: fill in the new Setpoint Voltage / Setpoint Current
#    sample values on row 2 (5 / 25).
# ---------------------------------------------------------------------------
$wsBmw = $wb.Worksheets.Item("bmw")
$wsBmw.Range("C2").Value = 5
$wsBmw.Range("D2").Value = 25

$wsDcb12 = $wb.Worksheets.Item("dcb1.2")
$wsDcb12.Range("C2").Value = 5
$wsDcb12.Range("D2").Value = 25

# ---------------------------------------------------------------------------
# 4) obc: add the Max Res value that was missing on row 2.
# ---------------------------------------------------------------------------
$wsObc = $wb.Worksheets.Item("obc")
$wsObc.Range("B2").Value = 10

# ---------------------------------------------------------------------------
# 5) Add the new "Setpoint Voltage" / "Setpoint Current" headers to every
#    sheet (columns C/D, row 1). Done last so these two brand-new shared
#    strings land at the very end of the shared-string table.
# ---------------------------------------------------------------------------
$sheetNames = @("bmw", "obc", "dcb1.2", "dcb1.2h", "dcb2.0", "5dh", "nissan", "custom")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C1").Value = "Setpoint Voltage"
    $ws.Range("D1").Value = "Setpoint Current"
    $ws.Columns.Item(3).ColumnWidth = 13.291666666666666
    $ws.Columns.Item(4).ColumnWidth = 13.291666666666666
}

# ---------------------------------------------------------------------------
# 6) Column width tweaks on column A (bestFit width drifted slightly once
#    the default font changed).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("bmw").Columns.Item(1).ColumnWidth = 10.666666666666666
$wb.Worksheets.Item("obc").Columns.Item(1).ColumnWidth = 10.666666666666666
$wb.Worksheets.Item("dcb1.2").Columns.Item(1).ColumnWidth = 14.541666666666666
$wb.Worksheets.Item("dcb1.2h").Columns.Item(1).ColumnWidth = 10.666666666666666
$wb.Worksheets.Item("dcb2.0").Columns.Item(1).ColumnWidth = 10.666666666666666
$wb.Worksheets.Item("5dh").Columns.Item(1).ColumnWidth = 10.666666666666666
$wb.Worksheets.Item("nissan").Columns.Item(1).ColumnWidth = 10.666666666666666
$wb.Worksheets.Item("custom").Columns.Item(1).ColumnWidth = 10.666666666666666

# ---------------------------------------------------------------------------
# 7) Selections (recorded cursor position per sheet) + re-activate bmw so it
#    stays the selected tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("obc").Activate()
$wb.Worksheets.Item("obc").Range("C2").Select()

$wb.Worksheets.Item("dcb1.2").Activate()
$wb.Worksheets.Item("dcb1.2").Range("C1:D2").Select()

$wb.Worksheets.Item("dcb1.2h").Activate()
$wb.Worksheets.Item("dcb1.2h").Range("C2:E2").Select()

$wb.Worksheets.Item("dcb2.0").Activate()
$wb.Worksheets.Item("dcb2.0").Range("C2:E3").Select()

$wb.Worksheets.Item("5dh").Activate()
$wb.Worksheets.Item("5dh").Range("C2:F2").Select()

$wb.Worksheets.Item("nissan").Activate()
$wb.Worksheets.Item("nissan").Range("C2:E2").Select()

$wb.Worksheets.Item("custom").Activate()
$wb.Worksheets.Item("custom").Range("C2:E2").Select()

$wb.Worksheets.Item("bmw").Activate()
$wb.Worksheets.Item("bmw").Range("D3").Select()
